# Weekly update: a new "Betarraga" price record (week of 2023-01-05) is
# inserted for "Terminal Hortofrutícola Agro Chillán" ahead of the existing
# history, pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 480; Excel shifts rows 480:510 down to 481:511 and the
# sheet dimension grows from A1:R510 to A1:R511 automatically.
$ws.Rows.Item(480).Insert()

# Populate the newly-opened row 480 with the new weekly record.
$ws.Range("A480").Value = 7
$ws.Range("B480").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C480").Value = "Ñuble"
$ws.Range("D480").Value = 44931
$ws.Range("E480").Value = 16
$ws.Range("F480").Value = 100114014
$ws.Range("G480").Value = "Betarraga"
$ws.Range("H480").Value = "Sin especificar"
$ws.Range("I480").Value = "Segunda"
$ws.Range("J480").Value = 200
$ws.Range("K480").Value = 600
$ws.Range("L480").Value = 600
$ws.Range("M480").Value = 600
$ws.Range("N480").Value = "`$/paquete 5 unidades"
$ws.Range("O480").Value = "Provincia de Diguillín"
$ws.Range("P480").Value = 120
$ws.Range("Q480").Value = 5
$ws.Range("R480").Value = "Hortaliza"
